$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first three data rows (2007年, 2008年, 2009年) -- everything
# shifts up by three rows so the former row 5 (2010年) becomes row 2, etc.
$ws.Range("A2:A4").EntireRow.Delete()

# Add a new trailing row for 2021年 (now row 13). Copy the formatting of the
# prior row (2020年) across the whole row first -- this both carries over the
# bold/centered/bordered style used on the column-A year labels and leaves
# the not-yet-reported columns (B, C, D, G, H, I, J) as blank cells, matching
# how the other "missing data" rows above are represented.
$ws.Range("A12:M12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "2021年"

$ws.Range("E13").Value = 185181
$ws.Range("F13").Value = 72791
$ws.Range("K13").Value = 29293
$ws.Range("L13").Value = 24951
$ws.Range("M13").Value = 54244
